$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-19 21:48:50'
$ws.Range("I2").Value = '3.7 mm'
$ws.Range("E3").Value = '2026-02-19 21:48:52'
$ws.Range("I3").Value = '5.6 mm'
$ws.Range("E4").Value = '2026-02-19 21:48:55'
$ws.Range("J4").Value = '1010.1 hPa'
$ws.Range("E5").Value = '2026-02-19 21:48:58'
$ws.Range("I5").Value = '7.9 mm'
$ws.Range("E6").Value = '2026-02-19 21:49:01'
$ws.Range("J6").Value = '1010.3 hPa'
$ws.Range("O6").Value = '10.3 °C'
$ws.Range("E7").Value = '2026-02-19 21:49:04'
$ws.Range("J7").Value = '1011.3 hPa'
$ws.Range("O7").Value = '13.8 °C'
$ws.Range("E8").Value = '2026-02-19 21:49:06'
$ws.Range("J8").Value = '1011.0 hPa'
$ws.Range("E9").Value = '2026-02-19 21:49:09'
$ws.Range("H9").Value = "'79%"
$ws.Range("O9").Value = '10.3 °C'
$ws.Range("E10").Value = '2026-02-19 21:49:12'
$ws.Range("N10").Value = '3.8 °C 21:02 TU'
$ws.Range("O10").Value = '10.1 °C'
$ws.Range("E11").Value = '2026-02-19 21:49:15'
$ws.Range("H11").Value = "'62%"
$ws.Range("O11").Value = '5.8 °C'
$ws.Range("E12").Value = '2026-02-19 21:49:17'
$ws.Range("E13").Value = '2026-02-19 21:49:20'
$ws.Range("H13").Value = "'61%"
$ws.Range("J13").Value = '1011.5 hPa'
$ws.Range("E14").Value = '2026-02-19 21:49:23'
$ws.Range("O14").Value = '13.2 °C'
$ws.Range("E15").Value = '2026-02-19 21:49:25'
$ws.Range("E16").Value = '2026-02-19 21:49:28'
$ws.Range("I16").Value = '10.7 mm'
$ws.Range("E17").Value = '2026-02-19 21:49:30'
$ws.Range("E18").Value = '2026-02-19 21:49:33'
$ws.Range("J18").Value = '1010.4 hPa'
$ws.Range("O18").Value = '11.5 °C'
$ws.Range("E19").Value = '2026-02-19 21:49:36'
$ws.Range("O19").Value = '5.2 °C'
$ws.Range("E20").Value = '2026-02-19 21:49:39'
$ws.Range("H20").Value = "'88%"
$ws.Range("E21").Value = '2026-02-19 21:49:42'
$ws.Range("J21").Value = '1011.6 hPa'
$ws.Range("O21").Value = '6.5 °C'
$ws.Range("E22").Value = '2026-02-19 21:49:44'
$ws.Range("I22").Value = '0.4 mm'
$ws.Range("L22").Value = '101.9 km/h - 327º 21:27 TU'
$ws.Range("E23").Value = '2026-02-19 21:49:47'
$ws.Range("I23").Value = '11.0 mm'
$ws.Range("E24").Value = '2026-02-19 21:49:50'
$ws.Range("J24").Value = '1015.2 hPa'
$ws.Range("O24").Value = '8.9 °C'
$ws.Range("E25").Value = '2026-02-19 21:49:53'
$ws.Range("H25").Value = "'68%"
$ws.Range("I25").Value = '7.1 mm'
$ws.Range("E26").Value = '2026-02-19 21:49:55'
$ws.Range("J26").Value = '1010.1 hPa'
$ws.Range("L26").Value = '71.6 km/h - 252º 21:07 TU'
$ws.Range("E27").Value = '2026-02-19 21:49:58'
$ws.Range("E28").Value = '2026-02-19 21:50:01'
$ws.Range("J28").Value = '1010.1 hPa'
$ws.Range("E29").Value = '2026-02-19 21:50:04'
$ws.Range("H29").Value = "'77%"
$ws.Range("N29").Value = '5.0 °C 21:17 TU'
$ws.Range("O29").Value = '10.3 °C'
$ws.Range("E30").Value = '2026-02-19 21:50:06'
$ws.Range("J30").Value = '1010.3 hPa'
$ws.Range("E31").Value = '2026-02-19 21:50:09'
$ws.Range("J31").Value = '1009.7 hPa'
$ws.Range("K31").Value = '8.6 MJ/m2'
$ws.Range("E32").Value = '2026-02-19 21:50:12'
$ws.Range("O32").Value = '4.7 °C'
$ws.Range("E33").Value = '2026-02-19 21:50:15'
$ws.Range("J33").Value = '1011.1 hPa'
$ws.Range("E34").Value = '2026-02-19 21:50:17'
$ws.Range("O34").Value = '-1.9 °C'
$ws.Range("E35").Value = '2026-02-19 21:50:20'
$ws.Range("J35").Value = '1016.6 hPa'
$ws.Range("E36").Value = '2026-02-19 21:50:22'
$ws.Range("J36").Value = '1010.5 hPa'
$ws.Range("E37").Value = '2026-02-19 21:50:25'
$ws.Range("J37").Value = '1011.6 hPa'
$ws.Range("O37").Value = '5.8 °C'
$ws.Range("E38").Value = '2026-02-19 21:50:28'
$ws.Range("E39").Value = '2026-02-19 21:50:30'
$ws.Range("O39").Value = '-6.0 °C'
$ws.Range("E40").Value = '2026-02-19 21:50:33'
$ws.Range("H40").Value = "'70%"
$ws.Range("J40").Value = '1012.8 hPa'
$ws.Range("E41").Value = '2026-02-19 21:50:36'
$ws.Range("H41").Value = "'39%"
$ws.Range("J41").Value = '1013.2 hPa'
$ws.Range("E42").Value = '2026-02-19 21:50:39'
$ws.Range("H42").Value = "'77%"
$ws.Range("O42").Value = '11.0 °C'
$ws.Range("E43").Value = '2026-02-19 21:50:41'
$ws.Range("O43").Value = '9.0 °C'
$ws.Range("E44").Value = '2026-02-19 21:50:44'
$ws.Range("I44").Value = '9.7 mm'
$ws.Range("O44").Value = '-6.1 °C'
$ws.Range("E45").Value = '2026-02-19 21:50:47'
$ws.Range("J45").Value = '1015.9 hPa'
$ws.Range("E46").Value = '2026-02-19 21:50:50'
$ws.Range("J46").Value = '1016.0 hPa'
